$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing last row (27) with the corrected/latest total_cases figure.
# C27 and D27 are formulas and will recalc automatically.
$ws.Range("B27").Value = 85435

# Append new row 28 (date 2020-03-27 -> serial 43917)
$ws.Range("A28").Value = 43917
$ws.Range("A28").NumberFormat = $ws.Range("A27").NumberFormat
$ws.Range("B28").Value = 104126
$ws.Range("C28").Formula = "=B28-B27"
$ws.Range("D28").Formula = "=C28/C27"

# Append new row 29 (date 2020-03-28 -> serial 43918)
$ws.Range("A29").Value = 43918
$ws.Range("A29").NumberFormat = $ws.Range("A27").NumberFormat
$ws.Range("B29").Value = 123578
$ws.Range("C29").Formula = "=B29-B28"
$ws.Range("D29").Formula = "=C29/C28"

# Mirror the cursor/selection position left behind in the authored workbook
$ws.Range("G29").Select()
